$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    # Force the literal text into the cell without Excel coercing
    # numeric-looking strings (e.g. "0.169", "0.0000177") into floats.
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "69.372.23"
$ws.Range("E2").Value = "  +2.75%  "
$ws.Range("D3").Value = "2.424.09"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue $ws "D5" "563.68"
$ws.Range("E5").Value = "  +2.12%  "
Set-TextValue $ws "D6" "166.28"
$ws.Range("E6").Value = "  +4.55%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.49%  "
Set-TextValue $ws "D9" "0.169"
$ws.Range("E9").Value = "  +6.94%  "
$ws.Range("D10").Value = "2.422.36"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("E12").Value = "  +1.97%  "
Set-TextValue $ws "D13" "4.66"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "69.381.99"
$ws.Range("E14").Value = "  +2.91%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws "D15" "0.0000177"
$ws.Range("E15").Value = "  +4.26%  "
$ws.Range("D16").Value = "2.870.29"
$ws.Range("E16").Value = "  -1.22%  "
Set-TextValue $ws "D17" "23.93"
$ws.Range("E17").Value = "  +4.27%  "
$ws.Range("D18").Value = "2.433.56"
$ws.Range("E18").Value = "  +0.40%  "
Set-TextValue $ws "D19" "10.79"
$ws.Range("E19").Value = "  +4.15%  "
Set-TextValue $ws "D20" "341.88"
$ws.Range("E20").Value = "  +4.04%  "
Set-TextValue $ws "D21" "7.14"
$ws.Range("E21").Value = "  +4.61%  "
Set-TextValue $ws "D22" "3.89"
$ws.Range("E22").Value = "  +2.88%  "
Set-TextValue $ws "D23" "1.96"
$ws.Range("E23").Value = "  +6.12%  "
$ws.Range("E24").Value = "  -0.07%  "
Set-TextValue $ws "D25" "66.06"
$ws.Range("E25").Value = "  +0.52%  "
Set-TextValue $ws "D26" "3.81"
$ws.Range("E26").Value = "  +5.50%  "
Set-TextValue $ws "D27" "8.51"
$ws.Range("E27").Value = "  +6.01%  "
$ws.Range("E28").Value = "  +0.25%  "
Set-TextValue $ws "D29" "0.999"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "0.0₃0847"
$ws.Range("E30").Value = "  +5.71%  "
Set-TextValue $ws "D31" "7.39"
$ws.Range("E31").Value = "  +5.32%  "
Set-TextValue $ws "D32" "1.23"
$ws.Range("E32").Value = "  +9.72%  "
Set-TextValue $ws "D33" "451.85"
$ws.Range("E33").Value = "  +8.49%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +1.44%  "
Set-TextValue $ws "D36" "158.64"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("E38").Value = "  +5.29%  "
$ws.Range("E39").Value = "  -0.01%  "
Set-TextValue $ws "D40" "18.20"
$ws.Range("E40").Value = "  +2.53%  "
Set-TextValue $ws "D41" "0.302"
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("E42").Value = "  +4.69%  "
Set-TextValue $ws "D43" "4.39"
$ws.Range("E43").Value = "  +3.84%  "
Set-TextValue $ws "D44" "37.85"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("E45").Value = "  +1.81%  "
Set-TextValue $ws "D46" "2.08"
$ws.Range("E46").Value = "  +5.05%  "
Set-TextValue $ws "D47" "134.83"
$ws.Range("E47").Value = "  +3.90%  "
Set-TextValue $ws "D48" "3.39"
$ws.Range("E48").Value = "  +2.77%  "
$ws.Range("E49").Value = "  +2.54%  "
Set-TextValue $ws "D50" "0.488"
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws "D51" "0.0936"
$ws.Range("E51").Value = "  +2.52%  "
